$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Abrechnungsperiode" text: 30.2.2024 (invalid date) -> 29.2.2024 ---
$ws.Range("E16").Value = "Abrechnungsperiode 1.2.2024 - 29.2.2024"

# --- DATUM value: was stored as literal text "30.2.2024" (invalid date),
#     now becomes the real date 29 Feb 2024 (leap year), keeping the
#     existing mm-dd-yy style already applied to H8 ---
$ws.Range("H8").Value = 45351

# --- New note row 6: footnote about different times on the 19th/21st ---
$ws.Range("I6").Font.Name = "Arial"
$ws.Range("I6").Font.Italic = $true
$ws.Range("I6").Font.Size = 9
$ws.Range("I6").Value = "* Am 19. und 21. sind die Zeiten unterschiedlich:"

# --- New note row 7: special outbound trip time ---
$ws.Range("I7").Font.Name = "Arial"
$ws.Range("I7").Font.Italic = $true
$ws.Range("I7").Font.Size = 9
$ws.Range("I7").Value = "Hinfahrt: Von Zuhause nach Dychrain, 9:00 Uhr"

# --- New note row 8: special return trip time (General number format,
#     not the date format inherited from the cell's previous content) ---
$ws.Range("I8").ClearFormats()
$ws.Range("I8").Font.Name = "Arial"
$ws.Range("I8").Font.Italic = $true
$ws.Range("I8").Font.Size = 9
$ws.Range("I8").Value = "Rückfahrt: Von Dychrain nach Hause, 12:00 Uhr "

# --- Rounded grand total (was 2180.2788, now rounded to cents) ---
$ws.Range("O42").Value = 2188.38

# --- VAT rate label update: 7.7% -> 8.1% ---
$ws.Range("L42").Value = "Summe Total inkl. 8.1% MwSt:"

# D22 holds "=O42" already; it recalculates automatically to 2188.38

# --- Update the on-screen selection to match the saved view state ---
$ws.Range("L43").Select()
